$wb = $excel.ActiveWorkbook

# --- "Basic Game rubric" sheet: update a few description cells ---
$rubric = $wb.Worksheets.Item("Basic Game rubric")

$rubric.Range("C3").Value = "Multiple enviromental hazards are animated. one hazard (a hand) has contact animations. Cannon that aims it barrel at the player as he moves. Boss that has a charge, jump and land animation. Upon taking damage the sprite flashes red. When the player / enemy dies, a procedural particle system is played."
$rubric.Range("C4").Value = "Almost all hazards can be killed / destroyed using the player gun that shoots projectiles which damage and kill said hazards. There are hazards that drop from the ceiling and damage the player on contact, hazards that rise from the ground whem the player gets close, cannons that aim and shoot at the player when in sight, a boss that jumps and shoots towards the player."
$rubric.Range("C5").Value = "At the moment there are more than 10 different gameobjects the player can interact with in different ways (shooting them, walking over/on them, walking past them, walking in sight)"
$rubric.Range("C6").Value = "Health and score UI implemented. Implemented following menus: main menu with link to options menu and quit button, options menu to adjust volume, pause menu that freezes the game and can restart / go back to menu. end of level screen to go to the next level/ back to menu. End of game screen that shows the saved highscore and current score, button to go back to menu."
$rubric.Range("C7").Value = "All actions have sound, there is a main soundtrack and a boss soundtrack when starting the boss encounter. The master volume can be adjusted using the up and down buttons and in the options menu."

# --- "Game extras" sheet: drop the no-longer-relevant rows, then fill in the ones that are now implemented ---
$extras = $wb.Worksheets.Item("Game extras")

$extras.Rows("8:16").Delete()

$extras.Range("A2").Value = "Load/save game system"
$extras.Range("B2").Value = 1
$extras.Range("C2").Value = "The highscore gets loaded from a save file and displayed at the end game screen. When the player beats the previous highscore, it gets saved to a txt file."

$extras.Range("A3").Value = "Load level from a file"
$extras.Range("B3").Value = 1
$extras.Range("C3").Value = "The level and hazards get loaded from an SVG file. all other objects (enemies, rising hands, falling spikes,..) get loaded from a custom txt file containing their properties like position, health and size."

$extras.Range("A4").Value = "A.I."
$extras.Range("B4").Value = 1
$extras.Range("C4").Value = "Simple canon AI that aims its barrel towards the player and shoots when they are in sight and in range. Boss AI that charges and jumps towards the player at an interval, boss has a barrel like the cannon that aims at the player and shoots with an interval."

$extras.Range("A5").Value = "Particle system"
$extras.Range("B5").Value = 1
$extras.Range("C5").Value = "Entirely procedural, object based particle system. Particles get pooled at compile time to save CPU usage and prevent particles from getting destroy / instanciated everytime a particle system is used."

$extras.Range("A6").Value = "Screenshake"
$extras.Range("B6").Value = 1
$extras.Range("C6").Value = "Added instead of parallax, as recommended by Tom."

$extras.Range("A7").Value = "Object pooling"
$extras.Range("B7").Value = 1
$extras.Range("C7").Value = "Object pooling is used to instanciate a pool of projectiles and particles at the beginning to prevent a projectile being instanciated / destoryed everytime the player shoots or a particle system plays"

# --- Active sheet / selection bookkeeping so the view matches the saved state ---
$rubric.Activate()
$rubric.Range("C7").Select()

$extras.Range("C13").Select()
